# Update "想去人数" (F column) counts on the sheets that list event data.
# Mapping of worksheet name -> { row -> new F value }
$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        2  = 15051
        3  = 19169
        5  = 143
        13 = 59
        14 = 170
        15 = 231
        17 = 1483
        22 = 8005
        24 = 33
        26 = 67
        27 = 1248
        29 = 6075
        30 = 117
        31 = 74
        32 = 172
        33 = 155
        34 = 292
        35 = 5473
        36 = 719
        37 = 18
    }
    "演出" = @{
        3 = 18
    }
    "全部类型" = @{
        2  = 15051
        3  = 19169
        5  = 143
        13 = 59
        14 = 170
        15 = 231
        17 = 1483
        23 = 8005
        25 = 33
        27 = 67
        28 = 1248
        30 = 18
        32 = 6075
        33 = 117
        34 = 74
        35 = 172
        36 = 155
        37 = 292
        38 = 5473
        39 = 719
        40 = 18
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowValues = $updates[$sheetName]
    foreach ($row in $rowValues.Keys) {
        $ws.Range("F$row").Value = $rowValues[$row]
    }
}
